$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: swap C1 and D1 text values
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "max"

# Row 2 data: swap C2 and D2, with D2 getting a new numeric value
$ws.Range("C2").Value = "s__Proteus mirabilis"
$ws.Range("D2").Value = 0.4880628629810283
